$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 79: update Hours (C79) from 3 to 3.25
$ws.Range("C79").Value = 3.25

# Update the Notes + Weekly Total (D79) text to include the new note
# (modifies existing shared string in place, keeping its original index)
$ws.Range("D79").Value = "Finish 4.8, 4.9, 4.10, 4.11, quiz, 1 small problem"

# Row 79: set Milestones (E79) to new note
# (adds a brand new shared string entry, appended after the one above)
$ws.Range("E79").Value = "Finish Lesson 4"

$wb.Save()
